$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap match data (keep index/date columns A-E fixed) for rows 8/9, 14/15, 38/39 ---
# Row 8
$ws.Range("F8").Value = "NagaWorld"
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = "Visakha"
$ws.Range("I8").Value = 4
$ws.Range("J8").Value = 2.75
$ws.Range("K8").Value = "11/08/2023 01:12"
$ws.Range("L8").Value = 2.52
$ws.Range("M8").Value = "12/08/2023 09:30"
$ws.Range("N8").Value = 3.53
$ws.Range("O8").Value = "11/08/2023 01:12"
$ws.Range("P8").Value = 3.64
$ws.Range("Q8").Value = "12/08/2023 11:02"
$ws.Range("R8").Value = 2.03
$ws.Range("S8").Value = "11/08/2023 01:12"
$ws.Range("T8").Value = 2.28
$ws.Range("U8").Value = "12/08/2023 09:30"
$ws.Range("V8").Value = "https://www.betexplorer.com/football/cambodia/cpl/nagaworld-visakha/xhYp6ed5/"

# Row 9
$ws.Range("F9").Value = "Dangkor"
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = "Prey Veng"
$ws.Range("I9").Value = 5
$ws.Range("J9").Value = 2.04
$ws.Range("K9").Value = "12/08/2023 05:12"
$ws.Range("L9").Value = 2.04
$ws.Range("M9").Value = "12/08/2023 12:07"
$ws.Range("N9").Value = 3.77
$ws.Range("O9").Value = "12/08/2023 05:12"
$ws.Range("P9").Value = 3.63
$ws.Range("Q9").Value = "12/08/2023 12:07"
$ws.Range("R9").Value = 2.73
$ws.Range("S9").Value = "12/08/2023 05:12"
$ws.Range("T9").Value = 2.94
$ws.Range("U9").Value = "12/08/2023 12:07"
$ws.Range("V9").Value = "https://www.betexplorer.com/football/cambodia/cpl/dangkor-senchey-prey-veng/b1Zt7ysa/"

# Row 14
$ws.Range("F14").Value = "Svay Rieng"
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = "NagaWorld"
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 1.65
$ws.Range("K14").Value = "19/08/2023 01:13"
$ws.Range("L14").Value = 1.41
$ws.Range("M14").Value = "20/08/2023 12:02"
$ws.Range("N14").Value = 3.62
$ws.Range("O14").Value = "19/08/2023 01:13"
$ws.Range("P14").Value = 4.7
$ws.Range("Q14").Value = "20/08/2023 12:15"
$ws.Range("R14").Value = 3.89
$ws.Range("S14").Value = "19/08/2023 01:13"
$ws.Range("T14").Value = 5.42
$ws.Range("U14").Value = "20/08/2023 12:15"
$ws.Range("V14").Value = "https://www.betexplorer.com/football/cambodia/cpl/svay-rieng-nagaworld/dd5BMGtn/"

# Row 15
$ws.Range("F15").Value = "Visakha"
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = "Dangkor"
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 1.13
$ws.Range("K15").Value = "20/08/2023 03:12"
$ws.Range("L15").Value = 1.11
$ws.Range("M15").Value = "20/08/2023 12:03"
$ws.Range("N15").Value = 7.66
$ws.Range("O15").Value = "20/08/2023 03:12"
$ws.Range("P15").Value = 8.59
$ws.Range("Q15").Value = "20/08/2023 12:04"
$ws.Range("R15").Value = 11.8
$ws.Range("S15").Value = "20/08/2023 03:12"
$ws.Range("T15").Value = 12.86
$ws.Range("U15").Value = "20/08/2023 12:04"
$ws.Range("V15").Value = "https://www.betexplorer.com/football/cambodia/cpl/visakha-dangkor-senchey/6XP12DeU/"

# Row 38
$ws.Range("F38").Value = "Svay Rieng"
$ws.Range("G38").Value = 2
$ws.Range("H38").Value = "Kirivong Sok Sen Chey"
$ws.Range("I38").Value = 1
$ws.Range("J38").Value = 1.11
$ws.Range("K38").Value = "21/10/2023 00:13"
$ws.Range("L38").Value = 1.25
$ws.Range("M38").Value = "22/10/2023 12:44"
$ws.Range("N38").Value = 7
$ws.Range("O38").Value = "21/10/2023 00:13"
$ws.Range("P38").Value = 5.69
$ws.Range("Q38").Value = "22/10/2023 12:51"
$ws.Range("R38").Value = 9.710000000000001
$ws.Range("S38").Value = "21/10/2023 00:13"
$ws.Range("T38").Value = 7.79
$ws.Range("U38").Value = "22/10/2023 12:46"
$ws.Range("V38").Value = "https://www.betexplorer.com/football/cambodia/cpl/svay-rieng-kirivong-sok-sen-chey/vaoqpBvA/"

# Row 39
$ws.Range("F39").Value = "Visakha"
$ws.Range("G39").Value = 2
$ws.Range("H39").Value = "Angkor Tiger"
$ws.Range("I39").Value = 1
$ws.Range("J39").Value = 1.15
$ws.Range("K39").Value = "21/10/2023 00:43"
$ws.Range("L39").Value = 1.19
$ws.Range("M39").Value = "22/10/2023 12:45"
$ws.Range("N39").Value = 6.29
$ws.Range("O39").Value = "21/10/2023 00:43"
$ws.Range("P39").Value = 6.51
$ws.Range("Q39").Value = "22/10/2023 12:45"
$ws.Range("R39").Value = 8.19
$ws.Range("S39").Value = "21/10/2023 00:43"
$ws.Range("T39").Value = 8.699999999999999
$ws.Range("U39").Value = "22/10/2023 12:45"
$ws.Range("V39").Value = "https://www.betexplorer.com/football/cambodia/cpl/visakha-angkor-tiger/8YkmqVgG/"

# --- Append new rows 50 and 51, copying number formats/styles from row 49 ---
$ws.Range("A49:V49").Copy() | Out-Null
$ws.Range("A50:V51").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 50
$ws.Range("A50").Value = 49
$ws.Range("B50").Value = "cambodia"
$ws.Range("C50").Value = "cpl"
$ws.Range("D50").Value = "2023-2024"
$ws.Range("E50").Value = 45235.40625
$ws.Range("F50").Value = "Angkor Tiger"
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = "Dangkor"
$ws.Range("I50").Value = 1
$ws.Range("J50").Value = 2.22
$ws.Range("K50").Value = "03/11/2023 23:13"
$ws.Range("L50").Value = 2.65
$ws.Range("M50").Value = "05/11/2023 09:44"
$ws.Range("N50").Value = 3.3
$ws.Range("O50").Value = "03/11/2023 23:13"
$ws.Range("P50").Value = 3.5
$ws.Range("Q50").Value = "05/11/2023 09:44"
$ws.Range("R50").Value = 2.54
$ws.Range("S50").Value = "03/11/2023 23:13"
$ws.Range("T50").Value = 2.14
$ws.Range("U50").Value = "05/11/2023 09:44"
$ws.Range("V50").Value = "https://www.betexplorer.com/football/cambodia/cpl/angkor-tiger-dangkor-senchey/hphaxKx1/"

# Row 51
$ws.Range("A51").Value = 50
$ws.Range("B51").Value = "cambodia"
$ws.Range("C51").Value = "cpl"
$ws.Range("D51").Value = "2023-2024"
$ws.Range("E51").Value = 45235.5
$ws.Range("F51").Value = "Svay Rieng"
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = "Tiffy Army"
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 1.44
$ws.Range("K51").Value = "04/11/2023 00:12"
$ws.Range("L51").Value = 1.93
$ws.Range("M51").Value = "05/11/2023 11:59"
$ws.Range("N51").Value = 4.2
$ws.Range("O51").Value = "04/11/2023 00:12"
$ws.Range("P51").Value = 4.09
$ws.Range("Q51").Value = "05/11/2023 11:59"
$ws.Range("R51").Value = 4.57
$ws.Range("S51").Value = "04/11/2023 00:12"
$ws.Range("T51").Value = 2.93
$ws.Range("U51").Value = "05/11/2023 11:59"
$ws.Range("V51").Value = "https://www.betexplorer.com/football/cambodia/cpl/svay-rieng-tiffy-army/pKlivt8l/"

"edit complete"